$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to Text format before assigning, so
# values like "7.42" or "0.998" are stored as text (matching the
# original inlineStr / shared-string text cells) instead of numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.679.79'
$ws.Range("E2").Value = '  -6.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.328.56'
$ws.Range("E3").Value = '  -2.62%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.28'
$ws.Range("E5").Value = '  -3.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.93'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.332.10'
$ws.Range("E8").Value = '  -2.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("E9").Value = '  -1.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.42'
$ws.Range("E10").Value = '  -2.57%  '

$ws.Range("E11").Value = '  -5.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -2.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.886.69'
$ws.Range("E13").Value = '  -3.17%  '

$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.313.73'
$ws.Range("E15").Value = '  -3.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("E16").Value = '  -5.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.34'
$ws.Range("E17").Value = '  -2.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '59.878.68'
$ws.Range("E18").Value = '  -5.87%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.35'
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.63'
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.02'
$ws.Range("E21").Value = '  -8.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '352.91'
$ws.Range("E22").Value = '  -8.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("E23").Value = '  -1.46%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.454.67'
$ws.Range("E25").Value = '  -2.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '68.84'
$ws.Range("E26").Value = '  -6.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000111'
$ws.Range("E27").Value = '  +0.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  +0.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  +5.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.51'
$ws.Range("E30").Value = '  +6.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.87'
$ws.Range("E31").Value = '  -0.95%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  -4.55%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.152'
$ws.Range("E33").Value = '  -2.08%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.357.19'
$ws.Range("E35").Value = '  -2.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.81'
$ws.Range("E36").Value = '  -0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.31'
$ws.Range("E37").Value = '  +2.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.80'
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '158.04'
$ws.Range("E40").Value = '  -3.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0754'
$ws.Range("E41").Value = '  -2.71%  '

$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.37'
$ws.Range("E43").Value = '  +0.82%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +7.33%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.81'
$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.745'
$ws.Range("E46").Value = '  -5.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.18'
$ws.Range("E47").Value = '  -0.67%  '

$ws.Range("E48").Value = '  -3.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.78'
$ws.Range("E49").Value = '  +0.83%  '

$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.43'
$ws.Range("E50").Value = '  +17.15%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.25'
$ws.Range("E51").Value = '  +9.80%  '
